$d = $word.ActiveDocument

# --- Change 1 (do this one first, while paragraph indices still match the
# pre-edit document): add a <w:lastRenderedPageBreak/> marker inside the
# run that carries the text "区域开始" in the third "区域开始:" paragraph
# (the one right after the "区域结束" / blank-paragraph pair that precedes
# the "等等，你不会指的是淋雨吧？" exchange). ---
$target = $d.Paragraphs.Item(46)
$full = $target.Range

$insertPoint = $d.Range($full.Start, $full.Start)
$runXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:r w:rsidRPr="003004AD">' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:hAnsi="宋体" w:cs="宋体" w:hint="eastAsia"/>' + `
      '<w:color w:val="70AD47" w:themeColor="accent6"/>' + `
    '</w:rPr>' + `
    '<w:lastRenderedPageBreak/>' + `
    '<w:t>区域开始</w:t>' + `
  '</w:r>' + `
'</w:p>'
$insertPoint.InsertXML($runXml)

# the original "区域开始" run got pushed right after the one we just
# inserted (InsertXML does not consume/replace - it splices in place),
# so remove that now-duplicate run.
$dupStart = $full.Start + 4
$dupEnd = $dupStart + 4
$dup = $d.Range($dupStart, $dupEnd)
$dup.Delete()

# --- Change 2: remove the empty centered paragraph that sits just before
# the "我:" paragraph (empty paragraph with jc=center, no run content).
# Done after Change 1 so the paragraph index used above is unaffected. ---
$emptyPara = $d.Paragraphs.Item(29)
$emptyPara.Range.Delete()
